# Auto-generated script to update "想去人数" (F column) values per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 244
$ws.Range("F5").Value = 9430
$ws.Range("F6").Value = 9430
$ws.Range("F7").Value = 603
$ws.Range("F8").Value = 102
$ws.Range("F9").Value = 158
$ws.Range("F10").Value = 281
$ws.Range("F12").Value = 152
$ws.Range("F13").Value = 183
$ws.Range("F14").Value = 436
$ws.Range("F15").Value = 12056
$ws.Range("F16").Value = 12056
$ws.Range("F24").Value = 238
$ws.Range("F25").Value = 39
$ws.Range("F28").Value = 159
$ws.Range("F29").Value = 2725
$ws.Range("F32").Value = 2099
$ws.Range("F33").Value = 68
$ws.Range("F37").Value = 1001
$ws.Range("F38").Value = 4193
$ws.Range("F39").Value = 3631
$ws.Range("F40").Value = 539
$ws.Range("F42").Value = 3055
$ws.Range("F43").Value = 1320
$ws.Range("F44").Value = 194
$ws.Range("F46").Value = 420
$ws.Range("F47").Value = 521
$ws.Range("F48").Value = 67
$ws.Range("F49").Value = 222
$ws.Range("F50").Value = 127
$ws.Range("F51").Value = 138

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 42
$ws.Range("F17").Value = 11
$ws.Range("F20").Value = 186
$ws.Range("F21").Value = 6
$ws.Range("F22").Value = 8
$ws.Range("F25").Value = 44

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 244
$ws.Range("F9").Value = 9430
$ws.Range("F10").Value = 603
$ws.Range("F11").Value = 102
$ws.Range("F12").Value = 158
$ws.Range("F13").Value = 281
$ws.Range("F15").Value = 152
$ws.Range("F16").Value = 183
$ws.Range("F17").Value = 12056
$ws.Range("F18").Value = 12056
$ws.Range("F25").Value = 39
$ws.Range("F29").Value = 159
$ws.Range("F30").Value = 2725
$ws.Range("F33").Value = 2099
$ws.Range("F34").Value = 68
$ws.Range("F39").Value = 1001
$ws.Range("F40").Value = 186
$ws.Range("F41").Value = 8
$ws.Range("F42").Value = 3631
$ws.Range("F43").Value = 3055
$ws.Range("F45").Value = 1320
$ws.Range("F46").Value = 194
$ws.Range("F47").Value = 420
$ws.Range("F49").Value = 521
$ws.Range("F50").Value = 67
$ws.Range("F51").Value = 222
